$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column K (column 11), shifting existing
# columns K.. right by 3. The new columns inherit formatting from column J.
$ws.Range("K1:M1").EntireColumn.Insert()

# Populate the new header cells with the new race categories, matching the
# style of the other header cells in that row (style index 3 -> same as
# neighbouring string headers).
$ws.Range("K1").Value = "Race Unknown"
$ws.Range("L1").Value = "Race Other"
$ws.Range("M1").Value = "Race Refused to Answer"
$ws.Range("K1:M1").Style = $ws.Range("N1").Style

# Reset the view back to the top-left and select L2 (matches the
# post-edit selection captured in the workbook).
$ws.Range("A1").Select()
$ws.Range("L2").Select()

# Workbook-level calculation option: switch reference style to R1C1.
$excel.ReferenceStyle = -4136
